# Reverse the order of the "Periodo Mora" / "Valor Mora" table (rows 16-61)
# on sheet "Hoja1": previous account-statement periods are removed and the
# list is rebuilt in the opposite (newest-first) order, as described by the
# commit message "Elimna EC anteriores y se agregan nuevos, se modifica base
# de datos".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rng = $ws.Range("E16:F61")
$vals = $rng.Value()

$rowCount = $vals.GetUpperBound(0)
$colCount = $vals.GetUpperBound(1)

$reversed = New-Object 'object[,]' $rowCount, $colCount
for ($i = 1; $i -le $rowCount; $i++) {
    for ($j = 1; $j -le $colCount; $j++) {
        $reversed[$i - 1, $j - 1] = $vals[$rowCount - $i + 1, $j]
    }
}

$rng.Value = $reversed
